$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13 (content rows for Objectives text + Docente responsavel)
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# Remove the stray column-A cell/style that Insert() copied down into the two new rows
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# Give the new B13/C13 and B14/C14 cells the same formatting as the other content cells
# (copy format only, so no new style entries get created in styles.xml)
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Now set the text for every affected cell (rows 13-27)
$ws.Range("B13").Value2 = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("C13").Value2 = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("B14").Value2 = "3444370 - Rita de Cássia Lacerda Brambilla Rodrigues"
$ws.Range("C14").Value2 = "3444370 - Rita de Cássia Lacerda Brambilla Rodrigues"
$ws.Range("A15").Value2 = "Programa resumido:"
$ws.Range("B15").Value2 = "Introdução. Conceitos básicos de instrumentação para controle de processos. Instrumentos de medição de pressão. Dispositivos de medição de nível. Instrumentos de medição de vazão. Dispositivos de medição de temperatura. Sensores online para medição de meios, células e gases. Introdução à teoria de controle. Controles aplicados a bioprocessos."
$ws.Range("C15").Value2 = "Introdução. Conceitos básicos de instrumentação para controle de processos. Instrumentos de medição de pressão. Dispositivos de medição de nível. Instrumentos de medição de vazão. Dispositivos de medição de temperatura. Sensores online para medição de meios, células e gases. Introdução à teoria de controle. Controles aplicados a bioprocessos."
$ws.Range("A16").Value2 = "Short syllabus:"
$ws.Range("B16").Value2 = "Introduction. Basic concepts of instrumentation for process control. Pressure measuring instruments. Level measuring devices. Flow measuring instruments. Temperature measuring devices. Online sensors for measurement of media, cells and gases. Introduction to control theory. Controls applied to bioprocesses"
$ws.Range("C16").Value2 = "Introduction. Basic concepts of instrumentation for process control. Pressure measuring instruments. Level measuring devices. Flow measuring instruments. Temperature measuring devices. Online sensors for measurement of media, cells and gases. Introduction to control theory. Controls applied to bioprocesses"
$ws.Range("A17").Value2 = "Programa:"
$ws.Range("B17").Value2 = "1.Introdução: aspectos gerais relativos à instrumentação e controle de bioprocessos.2.Conceitos básicos de instrumentação para controle de processos: definições dos elementos em uma malha de controle. Características gerais de instrumentos: classes de instrumentos e definições. Identificação e símbolos de instrumentos: padronização ISA, exemplos de simbologia.3.Instrumentos de medição de pressão: manômetros, diafragmas, cápsulas e foles, tubos de Bourdon e outros sensores de pressão.4.Dispositivos de medição de nível: medição direta de nível e medição indireta de nível.5.Instrumentos de medição de vazão: medidores deprimogênios, medidores lineares, medidores volumétricos e outros.6.Dispositivos de medição de temperatura: termômetros, termômetros com mola de pressão, dispositivos de temperatura de resistência, termistores, termopares e outros.7.Instrumentos de medição de pH, potencial redox, pressões parciais de oxigênio dissolvido e gás carbônico. Medição de potência de agitação e velocidade do impelidor. Sensores online para propriedades celulares e determinação da concentração total de biomassa. Determinação da concentração de biomassa ativa ou viável.8. Introdução à teoria de controle: principais problemas para o controle de bioprocessos. Definições básicas (controle manual, controle por realimentação - feedback, controle por antecipação - feedforward, ganho e atraso), componentes de um sistema de controle (sensor/transmissor e controlador e elementos finais de controle), ações de controle Liga-desliga (on-off), auto-operado, proporcional (P), proporcional-integral (PI), proporcional-derivativa (PD), proporcional-integral-derivativa (PDI). Introdução à interface de comunicação.9.Controles aplicados a bioprocessos: controle em malha aberta, controle por sistema regulatório, controle em cascata, controle por pré-alimentação e controle seguidor de trajetória e outros (controle ótimo, sistema adaptativo e sistema de controle por aprendizado)."
$ws.Range("C17").Value2 = "1.Introdução: aspectos gerais relativos à instrumentação e controle de bioprocessos.2.Conceitos básicos de instrumentação para controle de processos: definições dos elementos em uma malha de controle. Características gerais de instrumentos: classes de instrumentos e definições. Identificação e símbolos de instrumentos: padronização ISA, exemplos de simbologia.3.Instrumentos de medição de pressão: manômetros, diafragmas, cápsulas e foles, tubos de Bourdon e outros sensores de pressão.4.Dispositivos de medição de nível: medição direta de nível e medição indireta de nível.5.Instrumentos de medição de vazão: medidores deprimogênios, medidores lineares, medidores volumétricos e outros.6.Dispositivos de medição de temperatura: termômetros, termômetros com mola de pressão, dispositivos de temperatura de resistência, termistores, termopares e outros.7.Instrumentos de medição de pH, potencial redox, pressões parciais de oxigênio dissolvido e gás carbônico. Medição de potência de agitação e velocidade do impelidor. Sensores online para propriedades celulares e determinação da concentração total de biomassa. Determinação da concentração de biomassa ativa ou viável.8. Introdução à teoria de controle: principais problemas para o controle de bioprocessos. Definições básicas (controle manual, controle por realimentação - feedback, controle por antecipação - feedforward, ganho e atraso), componentes de um sistema de controle (sensor/transmissor e controlador e elementos finais de controle), ações de controle Liga-desliga (on-off), auto-operado, proporcional (P), proporcional-integral (PI), proporcional-derivativa (PD), proporcional-integral-derivativa (PDI). Introdução à interface de comunicação.9.Controles aplicados a bioprocessos: controle em malha aberta, controle por sistema regulatório, controle em cascata, controle por pré-alimentação e controle seguidor de trajetória e outros (controle ótimo, sistema adaptativo e sistema de controle por aprendizado)."
$ws.Range("A18").Value2 = "Syllabus:"
$ws.Range("B18").Value2 = "1.Introduction: general aspects related to the instrumentation and control of bioprocesses.2.Basic concepts of instrumentation for process control: definitions of the elements in a control mesh. General instrument characteristics: instrument classes and definitions. Instrument identification and symbols: ISA standardization, symbology examples.3.Pressure measuring instruments: manometers, diaphragms, capsule and bellows, Bourdon tubes and other pressure sensors.4.Level measuring devices: direct level measurement and indirect level measurement.5.Flow measurement instruments: pressure meters, linear meters, volumetric meters and others.6.Temperature measuring devices: thermometers, pressure spring thermometers, temperature resistance devices, thermistors, thermocouples and others.7.pH measuring instruments, redox potential, partial pressures of dissolved oxygen and carbon dioxide. Measurement of agitation power and impeller speed. Online sensors for cell properties and determination of total biomass concentration. Determination of active or viable biomass concentration.8.Introduction to control theory: main problems for the control of bioprocesses. Basic settings (manual control, feedback control, feedforward control, gain and delay), components of a control system (sensor / transmitter and controller and final control elements), control actions On/Off (P), proportional-integral (PI), proportional-derivative (PD), proportional-integral-derivative (PDI). Introduction to the communication interface.9.Controls applied to bioprocesses: open loop control, control by regulatory system, cascade control, pre-feed control and trajectory tracking control and others (optimal control, adaptive system and learning control system)."
$ws.Range("C18").Value2 = "1.Introduction: general aspects related to the instrumentation and control of bioprocesses.2.Basic concepts of instrumentation for process control: definitions of the elements in a control mesh. General instrument characteristics: instrument classes and definitions. Instrument identification and symbols: ISA standardization, symbology examples.3.Pressure measuring instruments: manometers, diaphragms, capsule and bellows, Bourdon tubes and other pressure sensors.4.Level measuring devices: direct level measurement and indirect level measurement.5.Flow measurement instruments: pressure meters, linear meters, volumetric meters and others.6.Temperature measuring devices: thermometers, pressure spring thermometers, temperature resistance devices, thermistors, thermocouples and others.7.pH measuring instruments, redox potential, partial pressures of dissolved oxygen and carbon dioxide. Measurement of agitation power and impeller speed. Online sensors for cell properties and determination of total biomass concentration. Determination of active or viable biomass concentration.8.Introduction to control theory: main problems for the control of bioprocesses. Basic settings (manual control, feedback control, feedforward control, gain and delay), components of a control system (sensor / transmitter and controller and final control elements), control actions On/Off (P), proportional-integral (PI), proportional-derivative (PD), proportional-integral-derivative (PDI). Introduction to the communication interface.9.Controls applied to bioprocesses: open loop control, control by regulatory system, cascade control, pre-feed control and trajectory tracking control and others (optimal control, adaptive system and learning control system)."
$ws.Range("A19").Value2 = "Avaliação:"
$ws.Range("A20").Value2 = "Método:"
$ws.Range("B20").Value2 = "Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2) e por um trabalho (T)."
$ws.Range("C20").Value2 = "Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2) e por um trabalho (T)."
$ws.Range("A21").Value2 = "Critério:"
$ws.Range("B21").Value2 = "Média do período normal = (P1 + P2 + T)/3"
$ws.Range("C21").Value2 = "Média do período normal = (P1 + P2 + T)/3"
$ws.Range("A22").Value2 = "Norma de recuperação:"
$ws.Range("B22").Value2 = "Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.
Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
$ws.Range("C22").Value2 = "Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.
Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
$ws.Range("A23").Value2 = "Bibliografia:"
$ws.Range("B23").Value2 = "Bibliografia
1. BAILEY, J.E., OLLIS, D.F. Biochemical Engineering Fundamental. 
New York: Mc-Graw Hill, 2nd edition,1986.
2. BU´LOCK, K. B. Biotecnologia Básica. Zaragoza: Editorial Acribia, 
1987.
3. SCHMIDELL, W., LIMA, U.A., AQUARONE, E., BORZANI, W. 
Biotecnologia Industrial - Engenharia Bioquímica (Vol 2), São Paulo: 
Edgard Blucher Ltda, 2001.
4. HARTNETT, T. Instrumentation and Control of Bioprocesses. In: LYDERSEN, 
B. K.; D'ELIA, N. A.; NELSON, K. L. Bioprocess Engineering: Systems, Equipment 
and Facilities. New York. John Wiley & Sons, Inc. 1994, p.413-468.
5. CINAR, A.; PARULEKAR, S.J.; ÜNDEY, C.; BIROL, G. Batch Fermentation  Modeling, Monitoring and Control. New York: Marcel Dekker, Inc., 2003.
6. SCHÜGERL, K. Measuring, Modeling and Control. In: REHM, H. J.; REED, G. Biotechnology. Vol 4. Weinheim: VCH, 2nd Edition, 1991."
$ws.Range("C23").Value2 = "Bibliografia
1. BAILEY, J.E., OLLIS, D.F. Biochemical Engineering Fundamental. 
New York: Mc-Graw Hill, 2nd edition,1986.
2. BU´LOCK, K. B. Biotecnologia Básica. Zaragoza: Editorial Acribia, 
1987.
3. SCHMIDELL, W., LIMA, U.A., AQUARONE, E., BORZANI, W. 
Biotecnologia Industrial - Engenharia Bioquímica (Vol 2), São Paulo: 
Edgard Blucher Ltda, 2001.
4. HARTNETT, T. Instrumentation and Control of Bioprocesses. In: LYDERSEN, 
B. K.; D'ELIA, N. A.; NELSON, K. L. Bioprocess Engineering: Systems, Equipment 
and Facilities. New York. John Wiley & Sons, Inc. 1994, p.413-468.
5. CINAR, A.; PARULEKAR, S.J.; ÜNDEY, C.; BIROL, G. Batch Fermentation  Modeling, Monitoring and Control. New York: Marcel Dekker, Inc., 2003.
6. SCHÜGERL, K. Measuring, Modeling and Control. In: REHM, H. J.; REED, G. Biotechnology. Vol 4. Weinheim: VCH, 2nd Edition, 1991."
$ws.Range("A24").Value2 = "Requisitos:"
$ws.Range("B25").Value2 = "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)
"
$ws.Range("C25").Value2 = "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)
"
$ws.Range("B26").Value2 = "LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)
"
$ws.Range("C26").Value2 = "LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)
"
$ws.Range("B27").Value2 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)
"
$ws.Range("C27").Value2 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)
"
